$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Range("H116").Value = 3960.2
$ws.Range("I116").Value = 4148
$ws.Range("J116").Value = 3835
$ws.Range("K116").Value = 4148
$ws.Range("L116").Value = 3835
$ws.Range("M116").Value = -706
$ws.Range("N116").Value = -10719

# Row 132
$ws.Range("H132").Value = 977.6539
$ws.Range("I132").Value = 894.9091
$ws.Range("J132").Value = 1432.75
$ws.Range("K132").Value = 2684.7273
$ws.Range("L132").Value = 4298.25
$ws.Range("M132").Value = -154.7273
$ws.Range("N132").Value = -9358.25

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 6844.909
$ws.Range("I45").Value = 6529.4
$ws.Range("K45").Value = 6529.4
$ws.Range("M45").Value = -6152.4

# Row 61
$ws.Range("H61").Value = 11200
$ws.Range("I61").Value = 9500
$ws.Range("J61").Value = 13750
$ws.Range("K61").Value = 9500
$ws.Range("L61").Value = 13750
$ws.Range("M61").Value = -9288
$ws.Range("N61").Value = -14174

# Row 122
$ws.Range("H122").Value = 2197.1428
$ws.Range("I122").Value = 2197.1428
$ws.Range("K122").Value = 6591.428400000001
$ws.Range("M122").Value = -4141.428400000001

# Row 132
$ws.Range("H132").Value = 4794.7666
$ws.Range("I132").Value = 2515.2173
$ws.Range("K132").Value = 7545.651899999999
$ws.Range("M132").Value = -5015.651899999999

# Row 136
$ws.Range("H136").Value = 11200
$ws.Range("I136").Value = 9500
$ws.Range("J136").Value = 13750
$ws.Range("K136").Value = 28500
$ws.Range("L136").Value = 41250
$ws.Range("M136").Value = -25950
$ws.Range("N136").Value = -46350

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3657.4666
$ws.Range("I20").Value = 3855.889
$ws.Range("K20").Value = 3855.889
$ws.Range("M20").Value = -3608.889

# Row 107
$ws.Range("H107").Value = 2921.0715
$ws.Range("I107").Value = 2543.3333
$ws.Range("K107").Value = 2543.3333
$ws.Range("M107").Value = -623.3332999999998

# Row 134
$ws.Range("H134").Value = 5575.048
$ws.Range("I134").Value = 3563.2942
$ws.Range("K134").Value = 10689.8826
$ws.Range("M134").Value = -8154.882599999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 45458736
$ws.Range("I31").Value = 111112280
$ws.Range("K31").Value = 111112280
$ws.Range("M31").Value = -111111985

# Row 34
$ws.Range("H34").Value = 45458736
$ws.Range("I34").Value = 111112280
$ws.Range("K34").Value = 111112280
$ws.Range("M34").Value = -111112078

# Row 58
$ws.Range("H58").Value = 4044.9473
$ws.Range("J58").Value = 6022.5884
$ws.Range("L58").Value = 6022.5884
$ws.Range("N58").Value = -6428.5884

# Row 107
$ws.Range("H107").Value = 1467.9678
$ws.Range("I107").Value = 1210.4348
$ws.Range("J107").Value = 2208.375
$ws.Range("K107").Value = 1210.4348
$ws.Range("L107").Value = 2208.375
$ws.Range("M107").Value = 709.5652
$ws.Range("N107").Value = -6048.375

# Row 135
$ws.Range("H135").Value = 99333.336
$ws.Range("J135").Value = 99333.336
$ws.Range("L135").Value = 99333.336
$ws.Range("N135").Value = -109473.336

# Row 136
$ws.Range("H136").Value = 4044.9473
$ws.Range("J136").Value = 6022.5884
$ws.Range("L136").Value = 18067.7652
$ws.Range("N136").Value = -23167.7652

$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 1216.2941
$ws.Range("I92").Value = 8002
$ws.Range("J92").Value = 792.1875
$ws.Range("K92").Value = 24006
$ws.Range("L92").Value = 2376.5625
$ws.Range("M92").Value = -22758
$ws.Range("N92").Value = -4872.5625

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 71806.64
$ws.Range("I2").Value = 430.91666
$ws.Range("K2").Value = 430.91666
$ws.Range("M2").Value = -317.91666

# Row 102
$ws.Range("H102").Value = 1322.4286
$ws.Range("I102").Value = 1003.4706
$ws.Range("K102").Value = 1003.4706
$ws.Range("M102").Value = 618.5294

# Row 113
$ws.Range("H113").Value = 52853.715
$ws.Range("I113").Value = 4997
$ws.Range("J113").Value = 88746.25
$ws.Range("K113").Value = 4997
$ws.Range("L113").Value = 88746.25
$ws.Range("M113").Value = -2827
$ws.Range("N113").Value = -93086.25

# Row 126
$ws.Range("H126").Value = 2951.3044
$ws.Range("I126").Value = 2419.1765
$ws.Range("J126").Value = 4459
$ws.Range("K126").Value = 7257.529500000001
$ws.Range("L126").Value = 13377
$ws.Range("M126").Value = -4787.529500000001
$ws.Range("N126").Value = -18317

# Row 132
$ws.Range("H132").Value = 7031.846
$ws.Range("I132").Value = 4448.8
$ws.Range("K132").Value = 13346.4
$ws.Range("M132").Value = -10816.4

# Row 133
$ws.Range("H133").Value = 110390
$ws.Range("J133").Value = 110390
$ws.Range("L133").Value = 110390
$ws.Range("N133").Value = -120510

$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 650
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 533.3333
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 533.3333
$ws.Range("M12").Value = -830
$ws.Range("N12").Value = -873.3333

# Row 61
$ws.Range("H61").Value = 7335.684
$ws.Range("I61").Value = 7326.5557
$ws.Range("J61").Value = 7500
$ws.Range("K61").Value = 7326.5557
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -7124.5557
$ws.Range("N61").Value = -7904

# Row 82
$ws.Range("H82").Value = 2411.1177
$ws.Range("I82").Value = 1458.9
$ws.Range("K82").Value = 1458.9
$ws.Range("M82").Value = -1097.9

# Row 85
$ws.Range("H85").Value = 2411.1177
$ws.Range("I85").Value = 1458.9
$ws.Range("K85").Value = 1458.9
$ws.Range("M85").Value = -210.9000000000001

# Row 100
$ws.Range("H100").Value = 41668084
$ws.Range("I100").Value = 62500650
$ws.Range("J100").Value = 2949
$ws.Range("K100").Value = 62500650
$ws.Range("L100").Value = 2949
$ws.Range("M100").Value = -62500109
$ws.Range("N100").Value = -4031

# Row 113
$ws.Range("H113").Value = 7335.684
$ws.Range("I113").Value = 7326.5557
$ws.Range("J113").Value = 7500
$ws.Range("K113").Value = 7326.5557
$ws.Range("L113").Value = 7500
$ws.Range("M113").Value = -5156.5557
$ws.Range("N113").Value = -11840

# Row 122
$ws.Range("H122").Value = 100006900
$ws.Range("I122").Value = 125004870
$ws.Range("K122").Value = 375014610
$ws.Range("M122").Value = -375012160

# Row 132
$ws.Range("H132").Value = 8248.538
$ws.Range("I132").Value = 8519.799999999999
$ws.Range("K132").Value = 25559.4
$ws.Range("M132").Value = -23029.4

# Row 136
$ws.Range("H136").Value = 4642.8887
$ws.Range("I136").Value = 2986.6
$ws.Range("K136").Value = 8959.799999999999
$ws.Range("M136").Value = -6409.799999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# Row 6
$ws.Range("H6").Value = 4431
$ws.Range("J6").Value = 4431
$ws.Range("L6").Value = 4431
$ws.Range("N6").Value = -4661

# Row 8
$ws.Range("H8").Value = 1000000
$ws.Range("I8").Value = 1000000
$ws.Range("K8").Value = 1000000
$ws.Range("M8").Value = -999860

# Row 132
$ws.Range("H132").Value = 3576.818
$ws.Range("I132").Value = 1233.25
$ws.Range("K132").Value = 3699.75
$ws.Range("M132").Value = -1169.75
